$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "cheat sheet" picture that was anchored on the worksheet
# (covered the sensitivity-analysis area). Deleting every shape mirrors
# the xdr:pic node being stripped out of drawing1.xml.
foreach ($shp in @($ws.Shapes)) {
    $shp.Delete()
}

# Bump the "Increase (Decrease) in Sales Volume" sensitivity input from
# 0% to 5% -- this single input drives every downstream formula in the
# pro-forma cash-flow model (rows 29-70) to recalculate.
$ws.Range("C6").Value = 0.05

# Match the reviewer's final scroll position / zoom / selection.
$win = $excel.ActiveWindow
$win.Zoom = 118
$ws.Range("G40").Select() | Out-Null
